# This script reproduces the "Updated cryptos list" GitHub Actions commit:
# it refreshes the Price (column D) and Volume(1h) (column E) figures for the
# cryptocurrency table, and re-orders the Kaspa / VeChain rows (40-41).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.334.30"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "3.215.41"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'608.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.06%  "
$ws.Range("D6").Value = "'156.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.215.94"
$ws.Range("E8").Value = "  +0.46%  "
$ws.Range("E9").Value = "  -1.80%  "
$ws.Range("D10").Value = "'0.160"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.66%  "
$ws.Range("D11").Value = "'5.70"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.66%  "
$ws.Range("E12").Value = "  -3.59%  "
$ws.Range("E13").Value = "  -0.60%  "
$ws.Range("D14").Value = "'38.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.18%  "
$ws.Range("D15").Value = "3.744.62"
$ws.Range("E15").Value = "  +0.42%  "
$ws.Range("D16").Value = "66.486.64"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("E17").Value = "  -3.27%  "
$ws.Range("D18").Value = "3.217.16"
$ws.Range("E18").Value = "  +0.49%  "
$ws.Range("E19").Value = "  +1.23%  "
$ws.Range("D20").Value = "'506.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.72%  "
$ws.Range("D21").Value = "'15.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.53%  "
$ws.Range("E22").Value = "  -1.53%  "
$ws.Range("D23").Value = "'8.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.70%  "
$ws.Range("D24").Value = "'14.59"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.89%  "
$ws.Range("D25").Value = "'85.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.85%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  -0.54%  "
$ws.Range("D28").Value = "'9.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.69%  "
$ws.Range("D29").Value = "'0.134"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +47.88%  "
$ws.Range("E30").Value = "  -1.04%  "
$ws.Range("D31").Value = "'6.95"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.19%  "
$ws.Range("D32").Value = "'2.89"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.05%  "
$ws.Range("D33").Value = "'28.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.46%  "
$ws.Range("E35").Value = "  -5.35%  "
$ws.Range("E36").Value = "  -2.04%  "
$ws.Range("D37").Value = "'501.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.77%  "
$ws.Range("D38").Value = "'55.37"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.88%  "
$ws.Range("D39").Value = "0.0₃0771"
$ws.Range("E39").Value = "  +12.54%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.0419"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.90%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'0.130"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.03%  "
$ws.Range("E42").Value = "  +4.99%  "
$ws.Range("D43").Value = "'8.72"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.25%  "
$ws.Range("E44").Value = "  -1.54%  "
$ws.Range("D45").Value = "2.924.76"
$ws.Range("E45").Value = "  +0.90%  "
$ws.Range("D46").Value = "'2.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'28.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.50%  "
$ws.Range("E48").Value = "  +1.70%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("E50").Value = "  -1.00%  "
$ws.Range("D51").Value = "'121.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.32%  "
